# Weekly price-sheet update: insert one new observation row above the
# existing data block (at row 244), pushing the previous rows 244-270
# down to 245-271, and populate the newly inserted row with the latest
# week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 244; Excel shifts rows 244:270 -> 245:271
# and extends the sheet's used range (dimension) accordingly.
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new weekly record.
$ws.Cells.Item(244, 1).Value = 6
$ws.Cells.Item(244, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(244, 3).Value = "Metropolitana"
$ws.Cells.Item(244, 4).Value = 44858
$ws.Cells.Item(244, 5).Value = 13
$ws.Cells.Item(244, 6).Value = 100112022
$ws.Cells.Item(244, 7).Value = "Arveja Verde"
$ws.Cells.Item(244, 8).Value = "Perfection"
$ws.Cells.Item(244, 9).Value = "Primera"
$ws.Cells.Item(244, 10).Value = 70
$ws.Cells.Item(244, 11).Value = 18000
$ws.Cells.Item(244, 12).Value = 20000
$ws.Cells.Item(244, 13).Value = 19000
$ws.Cells.Item(244, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(244, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(244, 16).Value = 760
$ws.Cells.Item(244, 17).Value = 25
$ws.Cells.Item(244, 18).Value = "Hortaliza"
